$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-29 Monday" "2025-09-30 Tuesday"

Replace-Text "60×80=" "27×69="
Replace-Text "76×85=" "20×97="
Replace-Text "42×29=" "78×87="
Replace-Text "26×55=" "62×67="
Replace-Text "44×90=" "88×19="

Replace-Text "48×39=" "80×17="
Replace-Text "45×52=" "79×74="
Replace-Text "40×78=" "21×29="
Replace-Text "87×97=" "21×85="
Replace-Text "44×43=" "19×33="

Replace-Text "87×44=" "35×14="
Replace-Text "87×99=" "51×97="
Replace-Text "35×67=" "40×92="
Replace-Text "69×19=" "20×29="
Replace-Text "39×43=" "55×34="

Replace-Text "43×86=" "81×12="
Replace-Text "99×30=" "50×72="
Replace-Text "96×45=" "51×25="
Replace-Text "65×78=" "75×52="
Replace-Text "83×87=" "88×26="

Replace-Text "38×36=" "53×17="
Replace-Text "19×64=" "27×51="
Replace-Text "43×43=" "27×99="
Replace-Text "68×61=" "60×11="
Replace-Text "69×78=" "79×40="
